# Week 13 logging update
# Appends new per-game stat samples to the YDS and ST sheets' running
# space-separated number lists, and updates the season-total cells on
# OFF, DEF, ST, TURNS and PEN accordingly.

$wb = $excel.ActiveWorkbook

function Append-Text {
    param($Range, [string]$Suffix)
    $current = $Range.Text
    $Range.Value = $current + $Suffix
}

function Add-Number {
    param($Range, [int]$Delta)
    $current = [int]$Range.Text
    $Range.Value = $current + $Delta
}

# ---------------------------------------------------------------------
# YDS sheet: append week 13 per-game yardage samples
# ---------------------------------------------------------------------
$ydsWs = $wb.Worksheets.Item("YDS")
Append-Text $ydsWs.Range("B2") " 1 12 5 18 -1 3 5 3 8 -2 4 5 3 8 6 8 0 -1 -2 13 3 2 2 3 1"
Append-Text $ydsWs.Range("C2") " 14 16 18 0 9 9 8 2 10 29 3 10 15 15 3 17 6 5 22 18 13 5 6"
Append-Text $ydsWs.Range("B3") " 2 7 0 5 -1 5 4 2 3 6 2 2 0 1 8 3 2 0 8 13 3 5 8 -2"
Append-Text $ydsWs.Range("C3") " 3 18 7 12 4 17 5 14 7 7 4 9 40 29 25 9 2 5 9 5 5"

# ---------------------------------------------------------------------
# OFF sheet: season totals through week 13
# ---------------------------------------------------------------------
$offWs = $wb.Worksheets.Item("OFF")
Add-Number $offWs.Range("C2") 14
Add-Number $offWs.Range("D2") 1
Add-Number $offWs.Range("E2") 1
Add-Number $offWs.Range("F2") 6
Add-Number $offWs.Range("G2") 1
Add-Number $offWs.Range("H2") 1
Add-Number $offWs.Range("J2") 1
Add-Number $offWs.Range("N2") 7
Add-Number $offWs.Range("B3") 1
Add-Number $offWs.Range("C3") 14
Add-Number $offWs.Range("E3") 6
Add-Number $offWs.Range("F3") 8
Add-Number $offWs.Range("G3") 1
Add-Number $offWs.Range("H3") 4
Add-Number $offWs.Range("I3") 8
Add-Number $offWs.Range("J3") 2
Add-Number $offWs.Range("L3") 37
Add-Number $offWs.Range("M3") 23
Add-Number $offWs.Range("Q3") 78

# ---------------------------------------------------------------------
# DEF sheet: season totals through week 13
# ---------------------------------------------------------------------
$defWs = $wb.Worksheets.Item("DEF")
Add-Number $defWs.Range("C2") 12
Add-Number $defWs.Range("F2") 4
Add-Number $defWs.Range("G2") 6
Add-Number $defWs.Range("J2") 2
Add-Number $defWs.Range("N2") 1
Add-Number $defWs.Range("C3") 12
Add-Number $defWs.Range("F3") 8
Add-Number $defWs.Range("G3") 2
Add-Number $defWs.Range("H3") 1
Add-Number $defWs.Range("I3") 4
Add-Number $defWs.Range("J3") 6
Add-Number $defWs.Range("L3") 32
Add-Number $defWs.Range("M3") 21
Add-Number $defWs.Range("Q3") 59

# ---------------------------------------------------------------------
# ST sheet: season totals & per-game samples through week 13
# ---------------------------------------------------------------------
$stWs = $wb.Worksheets.Item("ST")
Add-Number $stWs.Range("B2") 4
Add-Number $stWs.Range("D2") 4
Add-Number $stWs.Range("F2") 1
Add-Number $stWs.Range("G2") 1
Add-Number $stWs.Range("H2") 1
Add-Number $stWs.Range("J2") 2
Add-Number $stWs.Range("K2") 2
Add-Number $stWs.Range("B3") 1

Append-Text $stWs.Range("D3") " 55 45 46 57"
Append-Text $stWs.Range("B4") " 61 64 58"
Append-Text $stWs.Range("D4") " 14 5 6 8"
Append-Text $stWs.Range("B5") " 14 19 13"
Append-Text $stWs.Range("D5") " 0 0 12 0 0"
Append-Text $stWs.Range("B6") " 26 35 0"

# ---------------------------------------------------------------------
# TURNS sheet: season totals through week 13
# ---------------------------------------------------------------------
$turnsWs = $wb.Worksheets.Item("TURNS")
Add-Number $turnsWs.Range("B3") 1
Add-Number $turnsWs.Range("C3") 1
Add-Number $turnsWs.Range("D3") 1
Add-Number $turnsWs.Range("E3") 2

# ---------------------------------------------------------------------
# PEN sheet: season totals through week 13
# ---------------------------------------------------------------------
$penWs = $wb.Worksheets.Item("PEN")
Add-Number $penWs.Range("B2") 5
Add-Number $penWs.Range("D2") 1
Add-Number $penWs.Range("B3") 2
Add-Number $penWs.Range("D4") 1
